# Edit script: updates ZBP_03_strategie_domacnosti workbook with new data wave (22. 2. 2022)
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "data" sheet
$ws2 = $wb.Worksheets.Item(2)   # "pocetR" sheet

# ---- Sheet "data": fix AM1 typo year, add new AN column (22. 2. 2022) ----
$ws1.Range("AM1").Value = "25. 1. 2022"

# Give the new AN1 header the same look (bold/border/centered) as the other date headers
$ws1.Range("AL1").Copy()
$ws1.Range("AN1").PasteSpecial(-4122)
$ws1.Range("AN1").Value = "22. 2. 2022"

# AM (revised) / AN (new) data values for rows 2-45
$data1 = @{
    2 = @("0.22", "0.22")
    3 = @("0.12", "0.14")
    4 = @("0.5", "0.47")
    5 = @("0.29", "0.31")
    6 = @("0.17", "0.16")
    7 = @("0.22", "0.21")
    8 = @("0.27", "0.26")
    9 = @("0.22", "0.21")
    10 = @("0.22", "0.24")
    11 = @("0.22", "0.21")
    12 = @("0.23", "0.23")
    13 = @("0.38", "0.35")
    14 = @("0.21", "0.21")
    15 = @("0.22", "0.23")
    16 = @("0.24", "0.22")
    17 = @("0.19", "0.2")
    18 = @("0.28", "0.28")
    19 = @("0.29", "0.28")
    20 = @("0.19", "0.18")
    21 = @("0.15", "0.16")
    22 = @("0.18", "0.15")
    23 = @("0.23", "0.22")
    24 = @("0.44", "0.49")
    25 = @("0.36", "0.42")
    26 = @("0.12", "0.15")
    27 = @("0.09", "0.1")
    28 = @("0.13", "0.16")
    29 = @("0.16", "0.17")
    30 = @("0.09", "0.1")
    31 = @("0.1", "0.12")
    32 = @("0.13", "0.15")
    33 = @("0.14", "0.17")
    34 = @("0.2", "0.18")
    35 = @("0.09", "0.11")
    36 = @("0.14", "0.17")
    37 = @("0.12", "0.13")
    38 = @("0.06", "0.07000000000000001")
    39 = @("0.19", "0.27")
    40 = @("0.14", "0.16")
    41 = @("0.07000000000000001", "0.11")
    42 = @("0.09", "0.09")
    43 = @("0.07000000000000001", "0.1")
    44 = @("0.2", "0.16")
    45 = @("0.2", "0.28")
}
foreach ($r in $data1.Keys) {
    $pair = $data1[$r]
    $ws1.Cells.Item($r, 39).Value = [double]$pair[0]   # column AM
    $ws1.Cells.Item($r, 40).Value = [double]$pair[1]   # column AN
}

# Footer label date update (row 46)
$ws1.Range("A46").Value = "Život během pandemie, Strategie domácností, % respondentů celkově a ve skupinách, aktualizace 2. 3. 2022"

# ---- Sheet "pocetR": fix AL1 typo year, add new AM column (22. 2. 2022) ----
$ws2.Range("AL1").Value = "25. 1. 2022"

$ws2.Range("AK1").Copy()
$ws2.Range("AM1").PasteSpecial(-4122)
$ws2.Range("AM1").Value = "22. 2. 2022"

# AL (revised) / AM (new) sample-size values for rows 2-23
$data2 = @{
    2 = @(1848, 1786)
    3 = @(190, 170)
    4 = @(334, 348)
    5 = @(1324, 1268)
    6 = @(884, 860)
    7 = @(164, 159)
    8 = @(529, 511)
    9 = @(271, 256)
    10 = @(846, 820)
    11 = @(153, 146)
    12 = @(113, 109)
    13 = @(736, 711)
    14 = @(848, 827)
    15 = @(637, 610)
    16 = @(363, 349)
    17 = @(172, 173)
    18 = @(642, 588)
    19 = @(619, 623)
    20 = @(307, 297)
    21 = @(544, 504)
    22 = @(333, 316)
    23 = @(152, 181)
}
foreach ($r in $data2.Keys) {
    $pair = $data2[$r]
    $ws2.Cells.Item($r, 38).Value = $pair[0]   # column AL
    $ws2.Cells.Item($r, 39).Value = $pair[1]   # column AM
}

# Footer label date update (row 24) + extend the trailing blank-cell run into the new AM column
$ws2.Range("A24").Value = "Život během pandemie, Strategie domácností, velikost dotázaného souboru celkově a ve skupinách, aktualizace 2. 3. 2022"
$ws2.Range("AL24").Copy()
$ws2.Range("AM24").PasteSpecial(-4122)

